# Update "Industry to ISIC Code Map.xlsx": split the combined
# "ISIC 20T21: Chemicals and pharmaceutical products" column into two
# separate columns - "ISIC 20: Chemicals" and "ISIC 21: Pharmaceuticals" -
# on the ItICM sheet, and update the industry-to-ISIC flags accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ItICM")

$YELLOW = 65535
$NONE = -4142

# 1) Insert a new column before column K (the old combined
#    "ISIC 20T21" column). This shifts the old K column (and everything
#    to its right) one column to the right, preserving values/formatting.
$ws.Range("K1").EntireColumn.Insert()

# 2) Set the headers for the newly split columns.
#    New column K becomes "ISIC 20: Chemicals" (brand-new column, blank so far).
#    Column L (the old K column, shifted right) held the combined label and
#    must be renamed to "ISIC 21: Pharmaceuticals".
$ws.Range("K1").Value = "ISIC 20: Chemicals"
$ws.Range("K1").Style = $ws.Range("J1").Style
$ws.Range("L1").Value = "ISIC 21: Pharmaceuticals"

# 3) Fill in the new column K (ISIC 20: Chemicals) data values for rows 2-9.
#    Every row keeps a 0 except row 5 ("chemicals" industry), which was
#    previously flagged under the combined column and now flags specifically
#    under "ISIC 20: Chemicals".
$chemCol = @{2=0;3=0;4=0;5=1;6=0;7=0;8=0;9=0}
foreach ($r in 2..9) {
    $cell = $ws.Range("K" + $r)
    $cell.Value = $chemCol[$r]
    if ($chemCol[$r] -eq 1) {
        $cell.Interior.Color = $YELLOW
    } else {
        $cell.Interior.ColorIndex = $NONE
    }
}

# 4) Column L (ISIC 21: Pharmaceuticals) already carries over the old K
#    column's values via the insert/shift (row 5's old flag moved here).
#    Clear row 5 (chemicals no longer also flags pharmaceuticals) and add
#    the new flag on row 9 ("other industries" now also maps to the newly
#    separated pharmaceuticals category).
$ws.Range("L5").Value = 0
$ws.Range("L5").Interior.ColorIndex = $NONE

$ws.Range("L9").Value = 1
$ws.Range("L9").Interior.Color = $YELLOW
